# Applies the crypto price/volume refresh described by the commit diff.
# (Cronos/RenderToken rows 48-49 also swap places.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.374.67"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.942.63"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4633"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3869"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.93"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07825"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9755"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.62"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.30%  "
$ws.Range("D13").Value = "1.935.73"
$ws.Range("E13").Value = "  +3.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.077"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.757"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07053"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009809"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.10"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "29.408.23"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.470"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").Value = "2.165.80"
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.096"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.40"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.753"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.40"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.857"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09357"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8613"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.176"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.123"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05767"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.153"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.682"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5669"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1781"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.412"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.731"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +7.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002804"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +32.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5293"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.42"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.13%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.091"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.14%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06864"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.816"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.37"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.33%  "
